$d = $word.ActiveDocument

# Change 1: "2- Sistema seleciona..." -> "2- Sistema redireciona..."
$d.Content.Find.Execute(
    "2- Sistema seleciona usuário para página de suporte de acordo com a opção escolhida",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "2- Sistema redireciona usuário para página de suporte de acordo com a opção escolhida",
    2)

# Change 2: merge "8- Sistema não consegue enviar mensagem e exibe uma" with " de erro. "
$d.Content.Find.Execute(
    "8- Sistema não consegue enviar mensagem e exibe uma",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "8- Sistema não consegue enviar mensagem e exibe uma de erro.",
    2)

# Remove the now-redundant trailing run that used to hold " de erro. "
$d.Content.Find.Execute(
    " de erro. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    2)
